# Apply the "working with script 2" update to the Cell_Type_Markers_Summary table.
# The table grows from 22 data rows (A1:B22) to 27 data rows (A1:B27):
#  - "Crypt cells" is inserted before "Dendritic cells"
#  - "Fibroblasts" is inserted before "Foveolar cells"
#  - "Pericytes" is inserted before "Plasma cells"
#  - "Stromal cells" is inserted before "T cells"
#  - "Tuft cells" is appended after "T cells"
#  - several freq values are updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 5 new rows at their correct positions, shifting following rows down.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(27).Insert()

$data = @(
    @("B cells", 5),
    @("Basophils", 12),
    @("Crypt cells", 4),
    @("Dendritic cells", 14),
    @("Endothelial cells", 43),
    @("Enterochromaffin cells", 2),
    @("Enteroendocrine cells", 10),
    @("Eosinophils", 2),
    @("Fibroblasts", 14),
    @("Foveolar cells", 3),
    @("Gastric chief cells", 3),
    @("Macrophages", 11),
    @("Mast cells", 38),
    @("Metaplastic cells", 2),
    @("Monocytes", 3),
    @("Mucous neck cells", 1),
    @("Natural killer T cells", 5),
    @("Neutrophils", 3),
    @("Parietal cells", 2),
    @("Pericytes", 2),
    @("Plasma cells", 12),
    @("Plasmacytoid dendritic cells", 4),
    @("Proliferating Tff2+ cells", 3),
    @("Stromal cells", 1),
    @("T cells", 5),
    @("Tuft cells", 6)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

$ws.Range("A1:B27").Select()
